# Scheduled runner update: refresh market-price-derived profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns) across the
# per-job leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 75005
$ws.Range("J16").Value = 75005
$ws.Range("L16").Value = 75005
$ws.Range("N16").Value = -75465
$ws.Range("H40").Value = 1500.0785
$ws.Range("I40").Value = 1386.1515
$ws.Range("J40").Value = 1708.9445
$ws.Range("K40").Value = 1386.1515
$ws.Range("L40").Value = 1708.9445
$ws.Range("M40").Value = -1211.1515
$ws.Range("N40").Value = -2058.9445
$ws.Range("H53").Value = 815.9375
$ws.Range("I53").Value = 1386.75
$ws.Range("J53").Value = 245.125
$ws.Range("K53").Value = 1386.75
$ws.Range("L53").Value = 245.125
$ws.Range("M53").Value = -749.75
$ws.Range("N53").Value = -1519.125
$ws.Range("H64").Value = 3783.2727
$ws.Range("I64").Value = 3996.842
$ws.Range("J64").Value = 3493.4285
$ws.Range("K64").Value = 3996.842
$ws.Range("L64").Value = 3493.4285
$ws.Range("M64").Value = -3748.842
$ws.Range("N64").Value = -3989.4285
$ws.Range("H67").Value = 3783.2727
$ws.Range("I67").Value = 3996.842
$ws.Range("J67").Value = 3493.4285
$ws.Range("K67").Value = 3996.842
$ws.Range("L67").Value = 3493.4285
$ws.Range("M67").Value = -3138.842
$ws.Range("N67").Value = -5209.4285
$ws.Range("H92").Value = 38889830
$ws.Range("I92").Value = 7937352
$ws.Range("J92").Value = 111112280
$ws.Range("K92").Value = 7937352
$ws.Range("L92").Value = 111112280
$ws.Range("M92").Value = -7936104
$ws.Range("N92").Value = -111114776
$ws.Range("H106").Value = 133335020
$ws.Range("I106").Value = 41668530
$ws.Range("J106").Value = 500001000
$ws.Range("K106").Value = 41668530
$ws.Range("L106").Value = 500001000
$ws.Range("M106").Value = -41667899
$ws.Range("N106").Value = -500002262
$ws.Range("H113").Value = 1953.2727
$ws.Range("J113").Value = 1655.1428
$ws.Range("L113").Value = 1655.1428
$ws.Range("N113").Value = -8163.1428
$ws.Range("H116").Value = 7040.25
$ws.Range("I116").Value = 8764.643
$ws.Range("J116").Value = 3016.6667
$ws.Range("K116").Value = 8764.643
$ws.Range("L116").Value = 3016.6667
$ws.Range("M116").Value = -5322.643
$ws.Range("N116").Value = -9900.6667
$ws.Range("H137").Value = 1265.1666
$ws.Range("I137").Value = 1083.8518
$ws.Range("J137").Value = 1809.1111
$ws.Range("K137").Value = 3251.5554
$ws.Range("L137").Value = 5427.3333
$ws.Range("M137").Value = -701.5553999999997
$ws.Range("N137").Value = -10527.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4221.411
$ws.Range("I32").Value = 2969.016
$ws.Range("J32").Value = 11280.363
$ws.Range("K32").Value = 2969.016
$ws.Range("L32").Value = 11280.363
$ws.Range("M32").Value = -2682.016
$ws.Range("N32").Value = -11854.363
$ws.Range("H61").Value = 4176.5884
$ws.Range("I61").Value = 4360.125
$ws.Range("K61").Value = 4360.125
$ws.Range("M61").Value = -4148.125
$ws.Range("H63").Value = 111113500
$ws.Range("I63").Value = 125002390
$ws.Range("J63").Value = 2400
$ws.Range("K63").Value = 125002390
$ws.Range("L63").Value = 2400
$ws.Range("M63").Value = -125001704
$ws.Range("N63").Value = -3772
$ws.Range("H66").Value = 111113500
$ws.Range("I66").Value = 125002390
$ws.Range("J66").Value = 2400
$ws.Range("K66").Value = 625011950
$ws.Range("L66").Value = 12000
$ws.Range("M66").Value = -625008518
$ws.Range("N66").Value = -18864
$ws.Range("H74").Value = 4119.1816
$ws.Range("I74").Value = 1095.25
$ws.Range("J74").Value = 21053.2
$ws.Range("K74").Value = 1095.25
$ws.Range("L74").Value = 21053.2
$ws.Range("M74").Value = -221.25
$ws.Range("N74").Value = -22801.2
$ws.Range("H77").Value = 4119.1816
$ws.Range("I77").Value = 1095.25
$ws.Range("J77").Value = 21053.2
$ws.Range("K77").Value = 5476.25
$ws.Range("L77").Value = 105266
$ws.Range("M77").Value = -1108.25
$ws.Range("N77").Value = -114002
$ws.Range("H136").Value = 4176.5884
$ws.Range("I136").Value = 4360.125
$ws.Range("K136").Value = 13080.375
$ws.Range("M136").Value = -10530.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 997.4048
$ws.Range("I94").Value = 559.7273
$ws.Range("J94").Value = 2602.2222
$ws.Range("K94").Value = 559.7273
$ws.Range("L94").Value = 2602.2222
$ws.Range("M94").Value = -108.7273
$ws.Range("N94").Value = -3504.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5199.396
$ws.Range("I31").Value = 1911.125
$ws.Range("K31").Value = 1911.125
$ws.Range("M31").Value = -1616.125
$ws.Range("H34").Value = 5199.396
$ws.Range("I34").Value = 1911.125
$ws.Range("K34").Value = 1911.125
$ws.Range("M34").Value = -1709.125
$ws.Range("H105").Value = 1385.8572
$ws.Range("I105").Value = 1447.5
$ws.Range("K105").Value = 1447.5
$ws.Range("M105").Value = 299.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 533.3333
$ws.Range("I21").Value = 500
$ws.Range("J21").Value = 600
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 1800
$ws.Range("M21").Value = -1327
$ws.Range("N21").Value = -2146

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 45455692
$ws.Range("I113").Value = 100000790
$ws.Range("J113").Value = 1437.6666
$ws.Range("K113").Value = 100000790
$ws.Range("L113").Value = 1437.6666
$ws.Range("M113").Value = -99998620
$ws.Range("N113").Value = -5777.6666
$ws.Range("H126").Value = 6054.696
$ws.Range("I126").Value = 8246.532999999999
$ws.Range("J126").Value = 1945
$ws.Range("K126").Value = 24739.599
$ws.Range("L126").Value = 5835
$ws.Range("M126").Value = -22269.599
$ws.Range("N126").Value = -10775

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 71430720
$ws.Range("I40").Value = 100001810
$ws.Range("K40").Value = 100001810
$ws.Range("M40").Value = -100001674

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1600
$ws.Range("I81").Value = 1600
$ws.Range("K81").Value = 3200
$ws.Range("M81").Value = -2139
$ws.Range("H84").Value = 1600
$ws.Range("I84").Value = 1600
$ws.Range("K84").Value = 16000
$ws.Range("M84").Value = -10696
